$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I0, IF) in row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold/border/centered) from the existing H1 header
# cell onto the new header cells, reusing the same style record rather than
# creating a new one.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in data rows 2-18: I column is constant 1, J column mirrors H column
for ($r = 2; $r -le 18; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
